$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.854561805725098
$ws.Range("B1").Value = 2.065648555755615
$ws.Range("C1").Value = 2.079057455062866
$ws.Range("D1").Value = 1.829973220825195
$ws.Range("E1").Value = 1.363438367843628
